# Adjusted the runtime data: mirror the existing A:E table into a new
# I:M table with a second set of "running time (ms)" measurements, add
# AVERAGE formulas down column M (like column E already has), and add a
# trailing "Average" summary row (row 29) under both tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New J/K/L raw measurements for rows 3..28 (mirrors B/C/D) -----------
$data = @{
    3  = @(2.017, 2.997, 2.745)
    4  = @(0, 0, 0)
    5  = @(0.99, 1.017, 1.993)
    6  = @(2.996, 2.992, 4.003)
    7  = @(0.998, 0.999, 1)
    8  = @(0.975, 0, 0)
    9  = @(0.999, 1.996, 1.001)
    10 = @(4.986, 4, 5.004)
    11 = @(4.143, 4.011, 5)
    12 = @(0.996, 0.997, 0.994)
    13 = @(9.996, 11.001, 10.008)
    14 = @(0, 0, 0)
    15 = @(27.003, 27, 27.998)
    16 = @(1.005, 0.996, 0.997)
    17 = @(0.995, 0.995, 0.994)
    18 = @(0, 0, 0)
    19 = @(1, 0.996, 0.998)
    20 = @(0.996, 0.999, 1.001)
    21 = @(0, 0, 0)
    22 = @(0.94, 1.007, 1.001)
    23 = @(0, 0, 0)
    24 = @(0, 0, 0)
    25 = @(0.999, 0.998, 0.995)
    26 = @(15.001, 14.994, 14.992)
    27 = @(0, 0, 0)
    28 = @(0, 0, 0)
}

# --- Header block (rows 1-2), columns I:M mirroring A:E ------------------
# Row 1: merged label cell + merged "running time (ms)" title
$ws.Range("B1:E1").Copy()
$ws.Range("J1:M1").PasteSpecial(-4122)
$ws.Range("J1").Value = $ws.Range("B1").Text

$ws.Range("A1:A2").Copy()
$ws.Range("I1:I2").PasteSpecial(-4122)

# Row 2: t1 / t2 / t3 / average headers
$ws.Range("B2:E2").Copy()
$ws.Range("J2:M2").PasteSpecial(-4122)
$ws.Range("J2").Value = $ws.Range("B2").Text
$ws.Range("K2").Value = $ws.Range("C2").Text
$ws.Range("L2").Value = $ws.Range("D2").Text
$ws.Range("M2").Value = $ws.Range("E2").Text

# --- Data rows 3..28 -------------------------------------------------------
$ws.Range("A3:A28").Copy()
$ws.Range("I3:I28").PasteSpecial(-4122)

for ($r = 3; $r -le 28; $r++) {
    # Label column I mirrors column A (same letters A..Z)
    $ws.Cells.Item($r, 9).Value = $ws.Cells.Item($r, 1).Text

    $vals = $data[$r]
    $ws.Cells.Item($r, 10).Value = $vals[0]
    $ws.Cells.Item($r, 11).Value = $vals[1]
    $ws.Cells.Item($r, 12).Value = $vals[2]
    $ws.Cells.Item($r, 13).Formula = "=AVERAGE(J$r" + ":L$r)"
}

# --- New row 29: "Average" summary row for both tables ---------------------
$ws.Range("A28").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("I28").Copy()
$ws.Range("I29").PasteSpecial(-4122)

$ws.Range("A29").Value = "Average"
$ws.Range("I29").Value = "Average"
$ws.Range("E29").Formula = "=AVERAGE(E3:E28)"
$ws.Range("M29").Formula = "=AVERAGE(M3:M28)"

# --- Merge the new header cells (mirrors A1:A2 / B1:E1) --------------------
$ws.Range("I1:I2").Merge()
$ws.Range("J1:M1").Merge()

# --- Selection cosmetics (matches the saved workbook state) ----------------
[void]$ws.Range("F32").Select()

Write-Host "Runtime data adjusted."
